# Alloy vs UML - final update
# Splits the paragraph ending in "...formulas." into three paragraphs,
# inserting two new paragraphs of text and relocating the "_GoBack"
# bookmark plus a trailing space run to the end of the new, final
# paragraph.

$d = $word.ActiveDocument

$apos = [char]0x2019

# The "_GoBack" bookmark currently sits at the end of the "...formulas."
# paragraph. Remove it now; we'll re-create it in its new home once the
# new paragraphs/text exist.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Paragraph 6 is "There's also an association ... simplify some formulas."
$pFormulas = $d.Paragraphs.Item(6)

# --- New paragraph: "Point of interests ..." ---------------------------
$pFormulas.Range.InsertParagraphAfter()
$pPoint = $d.Paragraphs.Item(7)
$pPoint.Range.Text = "Point of interests are not modelled in alloy because they" + $apos + "re not useful in the model testing."

# --- New paragraph: "In the end, ..." -----------------------------------
$pPoint.Range.InsertParagraphAfter()
$pEnd = $d.Paragraphs.Item(8)
$pEnd.Range.Text = "In the end, some specific or multiple attributes are modelled only in alloy to avoid to full fill the UML diagram with less interesting aspects."

# Append the bookmark + a trailing space run at the end of this last
# paragraph (before its paragraph mark). Adding a bookmark collapsed
# directly at a paragraph-mark position is unreliable, so first place a
# throwaway marker character there, anchor the bookmark just before it,
# and then turn the marker into the real trailing space run.
$markPos = $pEnd.Range.End - 1
$filler = $d.Range($markPos, $markPos)
$filler.InsertBefore("Z")

$bmRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$fillerChar = $d.Range($markPos, $markPos + 1)
$fillerChar.Text = " "
